$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data, shifting rows 1-2 down to 2-3.
$null = $ws.Rows.Item(1).Insert()

# Add the new header labels in the freshly inserted row.
$ws.Range("B1").Value = "visits"
$ws.Range("C1").Value = "conversion"

# Update the selected cell to match the saved view state.
$null = $ws.Range("A7").Select()
